# Clean up vaccine price-list text:
#  - collapse embedded line breaks (Alt+Enter) inside cell text into a
#    single space so multi-line labels like "Fluzone\nQuadrivalent" read
#    as "Fluzone Quadrivalent" on one line
#  - strip the trailing footnote markers such as "[1]", "[2]", "[3]",
#    "[4]", "[5]" that were appended to vaccine names
#
# Applied across every worksheet in the workbook (Pediatric Vaccine,
# Adult Vaccine, Pediatric Influenza Vaccine, Adult Influenza Vaccine).

$wb = $excel.ActiveWorkbook
$newline = [char]10

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $used = $ws.UsedRange

    # Turn embedded line breaks into plain spaces.
    [void]$used.Replace($newline, " ")

    # Drop the footnote reference markers.
    [void]$used.Replace("[1]", "")
    [void]$used.Replace("[2]", "")
    [void]$used.Replace("[3]", "")
    [void]$used.Replace("[4]", "")
    [void]$used.Replace("[5]", "")
}
